$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text columns keep their text (inline-string-like) representation
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "68.112.18"
$ws.Range("E2").Value = "  +2.17%  "
$ws.Range("D3").Value = "3.341.24"
$ws.Range("E3").Value = "  +2.29%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "583.51"
$ws.Range("E5").Value = "  +2.64%  "
$ws.Range("D6").Value = "177.84"
$ws.Range("E6").Value = "  +1.71%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  +1.80%  "
$ws.Range("D9").Value = "3.340.23"
$ws.Range("E9").Value = "  +2.46%  "
$ws.Range("E10").Value = "  +5.57%  "
$ws.Range("E11").Value = "  +2.42%  "
$ws.Range("D12").Value = "47.01"
$ws.Range("E12").Value = "  +4.07%  "
$ws.Range("E13").Value = "  +2.50%  "
$ws.Range("D14").Value = "688.85"
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("D15").Value = "3.882.62"
$ws.Range("E15").Value = "  +2.55%  "
$ws.Range("D16").Value = "8.47"
$ws.Range("E16").Value = "  +2.34%  "
$ws.Range("D17").Value = "68.061.20"
$ws.Range("E17").Value = "  +2.09%  "
$ws.Range("D18").Value = "3.346.00"
$ws.Range("E18").Value = "  +2.47%  "
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("D21").Value = "11.09"
$ws.Range("E21").Value = "  +3.67%  "
$ws.Range("E22").Value = "  +1.72%  "
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").Value = "17.17"
$ws.Range("E23").Value = "  +1.87%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "5.36"
$ws.Range("E24").Value = "  +4.85%  "
$ws.Range("D25").Value = "98.66"
$ws.Range("E25").Value = "  +1.00%  "
$ws.Range("E26").Value = "  +1.30%  "
$ws.Range("E27").Value = "  +0.47%  "
$ws.Range("D28").Value = "9.57"
$ws.Range("E28").Value = "  +3.59%  "
$ws.Range("D29").Value = "33.09"
$ws.Range("E29").Value = "  +1.22%  "
$ws.Range("E30").Value = "  +2.78%  "
$ws.Range("D31").Value = "7.11"
$ws.Range("E31").Value = "  +6.37%  "
$ws.Range("D32").Value = "576.63"
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("D33").Value = "11.04"
$ws.Range("E33").Value = "  +2.85%  "
$ws.Range("E34").Value = "  +2.78%  "
$ws.Range("D35").Value = "3.722.04"
$ws.Range("E35").Value = "  -3.56%  "
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "57.00"
$ws.Range("E37").Value = "  +3.27%  "
$ws.Range("E38").Value = "  +3.00%  "
$ws.Range("D39").Value = "34.55"
$ws.Range("E39").Value = "  +9.43%  "
$ws.Range("E40").Value = "  +2.76%  "
$ws.Range("E41").Value = "  +3.14%  "
$ws.Range("D42").Value = "3.21"
$ws.Range("E42").Value = "  +5.99%  "
$ws.Range("D43").Value = "0.0₃0682"
$ws.Range("E43").Value = "  +2.83%  "
$ws.Range("D44").Value = "3.35"
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("E45").Value = "  +3.55%  "
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("E47").Value = "  +6.99%  "
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("E50").Value = "  -2.08%  "
$ws.Range("D51").Value = "129.40"
$ws.Range("E51").Value = "  +0.25%  "
